$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 217, shifting the existing rows 217-322 down to 218-323.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new record.
$ws.Cells.Item(217, 1).Value = 8
$ws.Cells.Item(217, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(217, 3).Value = "Coquimbo"
$ws.Cells.Item(217, 4).Value = 44839
$ws.Cells.Item(217, 5).Value = 4
$ws.Cells.Item(217, 6).Value = 100112012
$ws.Cells.Item(217, 7).Value = "Espinaca"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 2800
$ws.Cells.Item(217, 11).Value = 450
$ws.Cells.Item(217, 12).Value = 500
$ws.Cells.Item(217, 13).Value = 475
$ws.Cells.Item(217, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(217, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(217, 16).Value = 950
$ws.Cells.Item(217, 17).Value = 0.5
$ws.Cells.Item(217, 18).Value = "Hortaliza"
